$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Sort the data range A2:H13 by column H ("dissSet") ascending — this is
# the "Data > Sort" operation applied to the table (header row at row 1
# excluded from the sort, matching the workbook's existing sortCondition).
$dataRange = $ws.Range("A2:H13")
$keyRange = $ws.Range("H1")

$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($keyRange, 0, 1, $null, 0) | Out-Null
$sortObj.SetRange($dataRange)
$sortObj.Header = 2
$sortObj.MatchCase = $false
$sortObj.Orientation = 1
$sortObj.SortMethod = 1
$sortObj.Apply()

# Match the post-sort active selection recorded on Sheet1.
$ws.Range("D2:E13").Select()

$wb.Save()
